$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 15.438552
$ws.Range("H2").Value = 46.315656
$ws.Range("I2").Value = 0.4516083650784052
$ws.Range("J2").Value = 0.4516083650784052
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.557693
$ws.Range("N2").Value = 64.673079
$ws.Range("O2").Value = 0.5505707555812251
$ws.Range("P2").Value = 0.5505707555812251
$ws.Range("Q2").Value = 332.819564380536
$ws.Range("R2").Value = 2995.376079424824
$ws.Range("S2").Value = 0.2486423587880194
$ws.Range("T2").Value = 0.2486423587880193
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 15.438552
$ws.Range("H3").Value = 46.315656
$ws.Range("I3").Value = 0.4516083650784052
$ws.Range("J3").Value = 0.4516083650784052
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.981869333333333
$ws.Range("N3").Value = 11.945608
$ws.Range("O3").Value = 0.1016945926207894
$ws.Range("P3").Value = 0.1016945926207894
$ws.Range("Q3").Value = 61.474296759872
$ws.Range("R3").Value = 553.2686708388479
$ws.Range("S3").Value = 0.04592612871078914
$ws.Range("T3").Value = 0.04592612871078913
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 15.438552
$ws.Range("H4").Value = 46.315656
$ws.Range("I4").Value = 0.4516083650784052
$ws.Range("J4").Value = 0.4516083650784052
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.00696933333334
$ws.Range("N4").Value = 36.02090800000001
$ws.Range("O4").Value = 0.3066509100994217
$ws.Range("P4").Value = 0.3066509100994217
$ws.Range("Q4").Value = 185.370220415072
$ws.Range("R4").Value = 1668.331983735648
$ws.Range("S4").Value = 0.1384861161598049
$ws.Range("T4").Value = 0.1384861161598049
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 15.438552
$ws.Range("H5").Value = 46.315656
$ws.Range("I5").Value = 0.4516083650784052
$ws.Range("J5").Value = 0.4516083650784052
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.608641
$ws.Range("N5").Value = 4.825923
$ws.Range("O5").Value = 0.04108374169856382
$ws.Range("P5").Value = 0.04108374169856382
$ws.Range("Q5").Value = 24.835087727832
$ws.Range("R5").Value = 223.515789550488
$ws.Range("S5").Value = 0.01855376141979191
$ws.Range("T5").Value = 0.01855376141979191
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.632438333333333
$ws.Range("H6").Value = 4.897315
$ws.Range("I6").Value = 0.04775206941738987
$ws.Range("J6").Value = 0.04775206941738987
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 21.557693
$ws.Range("N6").Value = 64.673079
$ws.Range("O6").Value = 0.5505707555812251
$ws.Range("P6").Value = 0.5505707555812251
$ws.Range("Q6").Value = 35.19160443143166
$ws.Range("R6").Value = 316.724439882885
$ws.Range("S6").Value = 0.02629089293969946
$ws.Range("T6").Value = 0.02629089293969946
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.632438333333333
$ws.Range("H7").Value = 4.897315
$ws.Range("I7").Value = 0.04775206941738987
$ws.Range("J7").Value = 0.04775206941738987
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.981869333333333
$ws.Range("N7").Value = 11.945608
$ws.Range("O7").Value = 0.1016945926207894
$ws.Range("P7").Value = 0.1016945926207894
$ws.Range("Q7").Value = 6.500156138057777
$ws.Range("R7").Value = 58.50140524252
$ws.Range("S7").Value = 0.004856127246201118
$ws.Range("T7").Value = 0.004856127246201118
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.632438333333333
$ws.Range("H8").Value = 4.897315
$ws.Range("I8").Value = 0.04775206941738987
$ws.Range("J8").Value = 0.04775206941738987
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 12.00696933333334
$ws.Range("N8").Value = 36.02090800000001
$ws.Range("O8").Value = 0.3066509100994217
$ws.Range("P8").Value = 0.3066509100994217
$ws.Range("Q8").Value = 19.60063700689111
$ws.Range("R8").Value = 176.40573306202
$ws.Range("S8").Value = 0.01464321554597337
$ws.Range("T8").Value = 0.01464321554597337
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.632438333333333
$ws.Range("H9").Value = 4.897315
$ws.Range("I9").Value = 0.04775206941738987
$ws.Range("J9").Value = 0.04775206941738987
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.608641
$ws.Range("N9").Value = 4.825923
$ws.Range("O9").Value = 0.04108374169856382
$ws.Range("P9").Value = 0.04108374169856382
$ws.Range("Q9").Value = 2.626007232971666
$ws.Range("R9").Value = 23.634065096745
$ws.Range("S9").Value = 0.001961833685515935
$ws.Range("T9").Value = 0.001961833685515935
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 16.727748
$ws.Range("H10").Value = 50.183244
$ws.Range("I10").Value = 0.4893199132744809
$ws.Range("J10").Value = 0.4893199132744808
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 21.557693
$ws.Range("N10").Value = 64.673079
$ws.Range("O10").Value = 0.5505707555812251
$ws.Range("P10").Value = 0.5505707555812251
$ws.Range("Q10").Value = 360.611655965364
$ws.Range("R10").Value = 3245.504903688276
$ws.Range("S10").Value = 0.2694052343724705
$ws.Range("T10").Value = 0.2694052343724705
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 16.727748
$ws.Range("H11").Value = 50.183244
$ws.Range("I11").Value = 0.4893199132744809
$ws.Range("J11").Value = 0.4893199132744808
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.981869333333333
$ws.Range("N11").Value = 11.945608
$ws.Range("O11").Value = 0.1016945926207894
$ws.Range("P11").Value = 0.1016945926207894
$ws.Range("Q11").Value = 66.60770677692801
$ws.Range("R11").Value = 599.469360992352
$ws.Range("S11").Value = 0.04976118924168831
$ws.Range("T11").Value = 0.04976118924168831
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 16.727748
$ws.Range("H12").Value = 50.183244
$ws.Range("I12").Value = 0.4893199132744809
$ws.Range("J12").Value = 0.4893199132744808
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.00696933333334
$ws.Range("N12").Value = 36.02090800000001
$ws.Range("O12").Value = 0.3066509100994217
$ws.Range("P12").Value = 0.3066509100994217
$ws.Range("Q12").Value = 200.8495572517281
$ws.Range("R12").Value = 1807.646015265552
$ws.Range("S12").Value = 0.1500503967353897
$ws.Range("T12").Value = 0.1500503967353897
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 16.727748
$ws.Range("H13").Value = 50.183244
$ws.Range("I13").Value = 0.4893199132744809
$ws.Range("J13").Value = 0.4893199132744808
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.608641
$ws.Range("N13").Value = 4.825923
$ws.Range("O13").Value = 0.04108374169856382
$ws.Range("P13").Value = 0.04108374169856382
$ws.Range("Q13").Value = 26.908941270468
$ws.Range("R13").Value = 242.180471434212
$ws.Range("S13").Value = 0.02010309292493243
$ws.Range("T13").Value = 0.02010309292493242
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.3869703333333334
$ws.Range("H14").Value = 1.160911
$ws.Range("I14").Value = 0.01131965222972415
$ws.Range("J14").Value = 0.01131965222972415
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 21.557693
$ws.Range("N14").Value = 64.673079
$ws.Range("O14").Value = 0.5505707555812251
$ws.Range("P14").Value = 0.5505707555812251
$ws.Range("Q14").Value = 8.342187646107668
$ws.Range("R14").Value = 75.079688814969
$ws.Range("S14").Value = 0.006232269481035922
$ws.Range("T14").Value = 0.006232269481035922
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.3869703333333334
$ws.Range("H15").Value = 1.160911
$ws.Range("I15").Value = 0.01131965222972415
$ws.Range("J15").Value = 0.01131965222972415
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.981869333333333
$ws.Range("N15").Value = 11.945608
$ws.Range("O15").Value = 0.1016945926207894
$ws.Range("P15").Value = 0.1016945926207894
$ws.Range("Q15").Value = 1.540865303209778
$ws.Range("R15").Value = 13.867787728888
$ws.Range("S15").Value = 0.001151147422110807
$ws.Range("T15").Value = 0.001151147422110807
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.3869703333333334
$ws.Range("H16").Value = 1.160911
$ws.Range("I16").Value = 0.01131965222972415
$ws.Range("J16").Value = 0.01131965222972415
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.00696933333334
$ws.Range("N16").Value = 36.02090800000001
$ws.Range("O16").Value = 0.3066509100994217
$ws.Range("P16").Value = 0.3066509100994217
$ws.Range("Q16").Value = 4.646340925243113
$ws.Range("R16").Value = 41.81706832718801
$ws.Range("S16").Value = 0.003471181658253858
$ws.Range("T16").Value = 0.003471181658253858
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.3869703333333334
$ws.Range("H17").Value = 1.160911
$ws.Range("I17").Value = 0.01131965222972415
$ws.Range("J17").Value = 0.01131965222972415
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.608641
$ws.Range("N17").Value = 4.825923
$ws.Range("O17").Value = 0.04108374169856382
$ws.Range("P17").Value = 0.04108374169856382
$ws.Range("Q17").Value = 0.6224963439836666
$ws.Range("R17").Value = 5.602467095852999
$ws.Range("S17").Value = 0.0004650536683235589
$ws.Range("T17").Value = 0.0004650536683235589
